$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure purely-numeric-looking Price strings are written as text (matching
# the original inlineStr/text cell type) rather than being auto-converted to
# floating point numbers by Excel.

$ws.Range("D2").Value = '66.535.22'
$ws.Range("E2").Value = '  +0.50%  '
$ws.Range("D3").Value = '3.590.63'
$ws.Range("E3").Value = '  +0.76%  '
$ws.Range("E4").Value = '  -0.01%  '
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '608.81'
$ws.Range("E5").Value = '  +0.19%  '
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '148.36'
$ws.Range("E6").Value = '  +2.42%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -2.63%  '
$ws.Range("E9").Value = '  +1.10%  '
$ws.Range("E10").Value = '  -0.46%  '
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("D12").Value = '4.202.66'
$ws.Range("E12").Value = '  +0.84%  '
$ws.Range("E13").Value = '  +0.47%  '
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '29.80'
$ws.Range("E14").Value = '  -1.03%  '
$ws.Range("D15").Value = '3.614.30'
$ws.Range("E15").Value = '  +1.45%  '
$ws.Range("D16").Value = '66.656.19'
$ws.Range("E16").Value = '  +0.52%  '
$ws.Range("E17").Value = '  +0.83%  '
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = '11.49'
$ws.Range("E18").Value = '  +1.07%  '
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '6.37'
$ws.Range("E19").Value = '  +2.65%  '
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '15.05'
$ws.Range("E20").Value = '  +0.87%  '
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '426.96'
$ws.Range("E21").Value = '  -1.15%  '
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '0.616'
$ws.Range("E22").Value = '  +0.55%  '
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '78.79'
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("D24").Value = '3.741.83'
$ws.Range("E24").Value = '  +1.00%  '
$ws.Range("E25").Value = '  -0.02%  '
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '0.0000122'
$ws.Range("E26").Value = '  +3.29%  '
$ws.Range("E27").Value = '  +3.79%  '
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '9.37'
$ws.Range("E28").Value = '  +2.77%  '
$ws.Range("E29").Value = '  -0.48%  '
$ws.Range("E30").Value = '  +0.04%  '
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '1.48'
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("D32").Value = '3.587.79'
$ws.Range("E32").Value = '  +0.84%  '
$ws.Range("E33").Value = '  +2.90%  '
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '25.43'
$ws.Range("E34").Value = '  -0.63%  '
$ws.Range("E35").Value = '  -0.83%  '
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '5.63'
$ws.Range("E37").Value = '  +0.27%  '
$ws.Range("E38").Value = '  -2.44%  '
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '177.63'
$ws.Range("E39").Value = '  +3.13%  '
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '0.0855'
$ws.Range("E40").Value = '  -0.06%  '
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '5.25'
$ws.Range("E41").Value = '  +0.55%  '
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = '0.897'
$ws.Range("E42").Value = '  +0.11%  '
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '1.90'
$ws.Range("E43").Value = '  -0.83%  '
$ws.Range("E44").Value = '  +8.93%  '
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("E46").Value = '  -1.95%  '
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '24.95'
$ws.Range("E47").Value = '  -3.50%  '
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '23.93'
$ws.Range("E48").Value = '  +2.17%  '
$ws.Range("E49").Value = '  +0.57%  '
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '0.951'
$ws.Range("E50").Value = '  +0.59%  '
$ws.Range("D51").Value = '2.413.19'
$ws.Range("E51").Value = '  +4.53%  '
